$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  3/31/2025  Through  4/6/2025"

# --- Data table updates ---
$ws.Range("N14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = -100
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("I16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 5
$ws.Range("N14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 3
$ws.Range("I16").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = 5
$ws.Range("N14").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = -40
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -42.857142857142
$ws.Range("N16").Value = -87.096774193548
$ws.Range("I16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -68.75
$ws.Range("I17").Value = 19
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = -59.574468085106
$ws.Range("L17").Value = -13.636363636363
$ws.Range("M17").Value = 5.555555555555
$ws.Range("N17").Value = -66.666666666666
$ws.Range("I16").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 14
$ws.Range("K18").Value = -28.571428571428
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -23.076923076923
$ws.Range("N18").Value = -93.055555555555
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 10
$ws.Range("I19").Value = 27
$ws.Range("J19").Value = 37
$ws.Range("K19").Value = -27.027027027027
$ws.Range("L19").Value = -40
$ws.Range("M19").Value = 8
$ws.Range("N19").Value = -18.181818181818
$ws.Range("I16").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("I16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 2
$ws.Range("N14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = -70
$ws.Range("L20").Value = -14.285714285714
$ws.Range("M20").Value = 20
$ws.Range("N20").Value = -92.307692307692
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -50
$ws.Range("F21").Value = 26
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = -43.478260869565
$ws.Range("I21").Value = 71
$ws.Range("J21").Value = 133
$ws.Range("K21").Value = -46.616541353383
$ws.Range("L21").Value = -26.804123711340
$ws.Range("M21").Value = -7.792207792207
$ws.Range("N21").Value = -81.413612565445
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I16").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("I16").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -83.333333333333
$ws.Range("I23").Value = 4
$ws.Range("K23").Value = -60
$ws.Range("L23").Value = -33.333333333333
$ws.Range("M23").Value = -55.555555555555
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 33
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = -21.428571428571
$ws.Range("I24").Value = 98
$ws.Range("J24").Value = 114
$ws.Range("K24").Value = -14.035087719298
$ws.Range("L24").Value = -25.190839694656
$ws.Range("M24").Value = 42.028985507246
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 17.647058823529
$ws.Range("I25").Value = 53
$ws.Range("J25").Value = 47
$ws.Range("K25").Value = 12.765957446808
$ws.Range("L25").Value = -37.647058823529
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -66.666666666666
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 43
$ws.Range("J26").Value = 58
$ws.Range("K26").Value = -25.862068965517
$ws.Range("L26").Value = -14
$ws.Range("M26").Value = -40.277777777777
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("L27").Value = -50
